# Commit: "Update type of permit"
# Change the résumé's "Permit:" contact-details line from the Austrian
# "Red-White-Red Card plus" to "EU permanent residence".

$d = $word.ActiveDocument

# Locate the paragraph that holds the permit info so the replace is scoped
# precisely to it (the phrase is unique in the document, but this keeps the
# edit targeted and safe).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Permit:*") {
        $target = $p
    }
}

if ($target -ne $null) {
    $rng = $target.Range
    $rng.Find.Execute(
        "Red-White-Red Card plus",
        $true, $false, $false, $false, $false,
        $true, 1, $false,
        "EU permanent residence",
        2
    )
}
